$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A87").Value = "a"
$ws.Range("B87").Value = 0.0
$ws.Range("A88").Value = "a"
$ws.Range("B88").Value = 4168.0
